$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 2.14
$ws.Range("L2").Value = 1.29
$ws.Range("Q2").Value = 1.71
$ws.Range("S2").Value = 2.78
$ws.Range("T2").Value = 1.64
$ws.Range("X2").Value = 23
$ws.Range("AG2").Value = 15.5
$ws.Range("H3").Value = 3.4
$ws.Range("L3").Value = 1.29
$ws.Range("N3").Value = 4.6
$ws.Range("P3").Value = 2.24
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.74
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 2.38
$ws.Range("Z3").Value = 27
$ws.Range("AB3").Value = 12.5
$ws.Range("AN3").Value = 13
$ws.Range("F5").Value = 3.1
$ws.Range("I5").Value = 3.3
$ws.Range("F6").Value = 2.18
$ws.Range("G6").Value = 2.42
$ws.Range("P6").Value = 1.46
$ws.Range("Q6").Value = 2.72
$ws.Range("F10").Value = 2.92
$ws.Range("K10").Value = 2.96
$ws.Range("Q10").Value = 3.45
$ws.Range("U12").Value = 2.16
$ws.Range("G13").Value = 5
$ws.Range("J13").Value = 4.1
$ws.Range("N13").Value = 4.2
$ws.Range("Q13").Value = 1.82
$ws.Range("Z13").Value = 11.5
$ws.Range("AB13").Value = 18.5
$ws.Range("AC13").Value = 9.199999999999999
$ws.Range("H14").Value = 3.7
$ws.Range("I15").Value = 3.7
$ws.Range("U15").Value = 2.06
$ws.Range("J16").Value = 3.3
$ws.Range("K16").Value = 3.4
$ws.Range("R16").Value = 1.26
$ws.Range("AH16").Value = 21
$ws.Range("F17").Value = 3.4
$ws.Range("I17").Value = 2.68
$ws.Range("J17").Value = 2.84
$ws.Range("Q17").Value = 2.9
